$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying dummy dataset was regenerated with fewer rows: the table
# (and the sheet's used range) shrinks from A1:C62 down to A1:C38.
# Drop the now-unused trailing rows first ...
$ws.Range("A39:C62").EntireRow.Delete()

# ... then refresh rows 4-38 so they hold the new dataset's values.
$data = @(
    @(4,  "90 GOSSK St", "Smartphone", 990.33),
    @(5,  "90 GOSSK St", "Sofa",       983.65),
    @(6,  "54 BAZVE St", "Sofa",       105.84),
    @(7,  "54 BAZVE St", "Apple",      686.6),
    @(8,  "54 BAZVE St", "Smartphone", 434.85),
    @(9,  "54 BAZVE St", "Chair",      634.16),
    @(10, "54 BAZVE St", "Milk",       773.19),
    @(11, "39 EKKMV St", "Milk",       444.25),
    @(12, "39 EKKMV St", "Chair",      98.39),
    @(13, "97 SCZWD St", "Sofa",       71.64),
    @(14, "97 SCZWD St", "Apple",      129.41),
    @(15, "97 SCZWD St", "Chair",      420.1),
    @(16, "97 SCZWD St", "Laptop",     132.88),
    @(17, "97 SCZWD St", "Milk",       344.53),
    @(18, "60 TFERV St", "Chair",      711.38),
    @(19, "60 TFERV St", "Sofa",       184.24),
    @(20, "60 TFERV St", "Milk",       621.98),
    @(21, "60 TFERV St", "Smartphone", 548.14),
    @(22, "77 AROES St", "Milk",       804.93),
    @(23, "77 AROES St", "Apple",      892.47),
    @(24, "77 AROES St", "Smartphone", 719.05),
    @(25, "77 AROES St", "Laptop",     465.22),
    @(26, "29 OPTEP St", "Laptop",     400.53),
    @(27, "29 OPTEP St", "Smartphone", 677.01),
    @(28, "29 OPTEP St", "Milk",       360.96),
    @(29, "29 OPTEP St", "Chair",      156.53),
    @(30, "77 ILPWL St", "Milk",       138.53),
    @(31, "77 ILPWL St", "Chair",      553.65),
    @(32, "77 ILPWL St", "Smartphone", 929.13),
    @(33, "55 GVPCZ St", "Apple",      93.58),
    @(34, "55 GVPCZ St", "Sofa",       188.46),
    @(35, "55 GVPCZ St", "Smartphone", 873.02),
    @(36, "55 GVPCZ St", "Chair",      523.78),
    @(37, "55 GVPCZ St", "Milk",       341.67),
    @(38, "42 YBSAC St", "Chair",      440.05)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
}

# Match the new active selection recorded for the refreshed sheet.
$ws.Range("B6").Select()
